# Insert three new paragraphs ("Character", "Note:", and the explanatory
# sentence about Clyde's color palette) right after the paragraph that
# contains "Color Palette Options:" and before the paragraph holding the
# forest-path image.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Color Palette Options:*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find paragraph containing 'Color Palette Options:'"
}

# --- Paragraph: "Character" ---------------------------------------------
$target.Range.InsertParagraphAfter()
$charPara = $target.Next()
$charPara.Range.InsertAfter("Character")

# --- Paragraph: "Note:" --------------------------------------------------
$charPara.Range.InsertParagraphAfter()
$notePara = $charPara.Next()
$notePara.Range.InsertAfter("Note:")

# --- Paragraph: the descriptive sentence (no underline) ------------------
$notePara.Range.InsertParagraphAfter()
$descPara = $notePara.Next()
$descPara.Range.InsertAfter([string]::Concat("Clyde", [char]0x2019, "s color palette is lighter than that of the background.  The background will be darker colors."))
$descPara.Range.Underline = 0
